# Auto-generated edit script applying cryptos list update
# (commit: "Updated cryptos list on Sun Apr  2 12:21:09 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.472.77"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "'1.825.40"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'316.73"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").Value = "'0.5137"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "'0.3857"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "'0.08232"
$ws.Range("E9").Value = "  +7.83%  "
$ws.Range("D10").Value = "'1.122"
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D11").Value = "'41.95"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "'6.375"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").Value = "'21.08"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "'1.005"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").Value = "'7.471"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "'1.818.89"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "'94.16"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "'0.00001117"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("D19").Value = "'0.06628"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "'17.81"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "'6.053"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").Value = "'28.503.09"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "'11.56"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'159.72"
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'20.99"
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("D28").Value = "'2.031.76"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'2.409"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").Value = "'125.66"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").Value = "'0.1101"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").Value = "'1.098"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").Value = "'5.732"
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").Value = "'0.07475"
$ws.Range("E34").Value = "  +6.16%  "
$ws.Range("D35").Value = "'3.685"
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("D36").Value = "'12.49"
$ws.Range("E36").Value = "  +11.26%  "
$ws.Range("D37").Value = "'0.2220"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'0.02361"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("D39").Value = "'5.225"
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("D40").Value = "'8.837"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("D42").Value = "'1.186"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "'1.391"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'13.67"
$ws.Range("E44").Value = "  +2.11%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6197"
$ws.Range("E45").Value = "  +4.98%  "
$ws.Range("D46").Value = "'3.807"
$ws.Range("E46").Value = "  +2.59%  "
$ws.Range("D47").Value = "'128.02"
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("D48").Value = "'2.017"
$ws.Range("D49").Value = "'1.206"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "'0.06932"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "'1.079"
$ws.Range("E51").Value = "  +1.36%  "

Write-Host "Applied 103 cell updates to cryptos sheet"
